$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 36, shifting existing rows 36:51 down to 37:52
$ws.Rows.Item(36).Insert()

# Populate the new row 36 with the new weekly record
$ws.Cells.Item(36, 1).Value = 4
$ws.Cells.Item(36, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(36, 3).Value = "Los Lagos"
$ws.Cells.Item(36, 4).Value = 44855
$ws.Cells.Item(36, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(36, 5).Value = 10
$ws.Cells.Item(36, 6).Value = 300000000
$ws.Cells.Item(36, 7).Value = "Espárragos"
$ws.Cells.Item(36, 8).Value = "Sin especificar"
$ws.Cells.Item(36, 9).Value = "Primera"
$ws.Cells.Item(36, 10).Value = 600
$ws.Cells.Item(36, 11).Value = 1600
$ws.Cells.Item(36, 12).Value = 1700
$ws.Cells.Item(36, 13).Value = 1650
$ws.Cells.Item(36, 14).Value = "$/kilo"
$ws.Cells.Item(36, 15).Value = "Provincia de Linares"
$ws.Cells.Item(36, 16).Value = 1650
$ws.Cells.Item(36, 17).Value = 1
$ws.Cells.Item(36, 18).Value = "Hortaliza"
